# The workbook tracks weekly "Haba" (fava bean) price observations for the
# "Macroferia Regional de Talca" market. A new weekly observation is being
# inserted as the new row 16, pushing the previously-existing rows 16-20
# down to rows 17-21 (their data is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 16; this shifts the old
# rows 16-20 down to 17-21 automatically, carrying their values/styles.
$ws.Rows(16).Insert()

# Populate the new row 16 with the new weekly data point.
$ws.Cells.Item(16, 1).Value  = 5
$ws.Cells.Item(16, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(16, 3).Value  = "Maule"
$ws.Cells.Item(16, 4).Value  = 44455
$ws.Cells.Item(16, 5).Value  = 7
$ws.Cells.Item(16, 6).Value  = 100112026
$ws.Cells.Item(16, 7).Value  = "Haba"
$ws.Cells.Item(16, 8).Value  = "Sin especificar"
$ws.Cells.Item(16, 9).Value  = "Primera"
$ws.Cells.Item(16, 10).Value = 300
$ws.Cells.Item(16, 11).Value = 12000
$ws.Cells.Item(16, 12).Value = 12000
$ws.Cells.Item(16, 13).Value = 12000
$ws.Cells.Item(16, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(16, 15).Value = "Región del Maule"
$ws.Cells.Item(16, 16).Value = 480
$ws.Cells.Item(16, 17).Value = 25
$ws.Cells.Item(16, 18).Value = "Hortaliza"
